$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 194.35715
$ws.Range("I11").Value = 194.35715
$ws.Range("K11").Value = 194.35715
$ws.Range("M11").Value = -54.35714999999999
$ws.Range("H68").Value = 46666.332
$ws.Range("J68").Value = 54999.5
$ws.Range("L68").Value = 54999.5
$ws.Range("N68").Value = -56497.5
$ws.Range("H71").Value = 46666.332
$ws.Range("J71").Value = 54999.5
$ws.Range("L71").Value = 164998.5
$ws.Range("N71").Value = -172486.5
$ws.Range("H74").Value = 9059.3125
$ws.Range("I74").Value = 8925.714
$ws.Range("K74").Value = 8925.714
$ws.Range("M74").Value = -7989.714
$ws.Range("H77").Value = 9059.3125
$ws.Range("I77").Value = 8925.714
$ws.Range("K77").Value = 44628.57
$ws.Range("M77").Value = -39948.57
$ws.Range("H132").Value = 19353.125
$ws.Range("I132").Value = 21732.344
$ws.Range("J132").Value = 2698.6
$ws.Range("K132").Value = 65197.03200000001
$ws.Range("L132").Value = 8095.799999999999
$ws.Range("M132").Value = -62667.03200000001
$ws.Range("N132").Value = -13155.8
$ws.Range("H135").Value = 2693.9473
$ws.Range("I135").Value = 2121.923
$ws.Range("J135").Value = 3933.3333
$ws.Range("K135").Value = 19097.307
$ws.Range("L135").Value = 35399.9997
$ws.Range("M135").Value = -16562.307
$ws.Range("N135").Value = -40469.9997
$ws.Range("H138").Value = 16806.717
$ws.Range("I138").Value = 1277.2128
$ws.Range("K138").Value = 3831.6384
$ws.Range("M138").Value = 1308.3616
$ws.Range("H141").Value = 700
$ws.Range("I141").Value = 700
$ws.Range("K141").Value = 2100
$ws.Range("M141").Value = 3080

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19363.143
$ws.Range("I32").Value = 20034.018
$ws.Range("K32").Value = 20034.018
$ws.Range("M32").Value = -19747.018
$ws.Range("H45").Value = 3488.3125
$ws.Range("I45").Value = 1865.625
$ws.Range("K45").Value = 1865.625
$ws.Range("M45").Value = -1488.625
$ws.Range("H61").Value = 6084.5
$ws.Range("I61").Value = 955.9375
$ws.Range("K61").Value = 955.9375
$ws.Range("M61").Value = -743.9375
$ws.Range("H136").Value = 6084.5
$ws.Range("I136").Value = 955.9375
$ws.Range("K136").Value = 2867.8125
$ws.Range("M136").Value = -317.8125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 125000
$ws.Range("J50").Value = 125000
$ws.Range("L50").Value = 125000
$ws.Range("N50").Value = -126148
$ws.Range("H134").Value = 9443.788
$ws.Range("I134").Value = 10357.444
$ws.Range("K134").Value = 31072.332
$ws.Range("M134").Value = -28537.332

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2470
$ws.Range("I16").Value = 1261
$ws.Range("K16").Value = 1261
$ws.Range("M16").Value = -974
$ws.Range("H62").Value = 7112.6924
$ws.Range("I62").Value = 9311.875
$ws.Range("J62").Value = 3594
$ws.Range("K62").Value = 9311.875
$ws.Range("L62").Value = 3594
$ws.Range("M62").Value = -8687.875
$ws.Range("N62").Value = -4842
$ws.Range("H65").Value = 7112.6924
$ws.Range("I65").Value = 9311.875
$ws.Range("J65").Value = 3594
$ws.Range("K65").Value = 46559.375
$ws.Range("L65").Value = 17970
$ws.Range("M65").Value = -43439.375
$ws.Range("N65").Value = -24210
$ws.Range("H86").Value = 37592.72
$ws.Range("I86").Value = 57640
$ws.Range("K86").Value = 57640
$ws.Range("M86").Value = -56517
$ws.Range("H89").Value = 37592.72
$ws.Range("I89").Value = 57640
$ws.Range("K89").Value = 288200
$ws.Range("M89").Value = -282584
$ws.Range("H113").Value = 2470
$ws.Range("I113").Value = 1261
$ws.Range("K113").Value = 1261
$ws.Range("M113").Value = 909
$ws.Range("H122").Value = 2583.182
$ws.Range("J122").Value = 2957
$ws.Range("L122").Value = 8871
$ws.Range("N122").Value = -13771

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 75.71429000000001
$ws.Range("I26").Value = 60
$ws.Range("K26").Value = 180
$ws.Range("M26").Value = 108
$ws.Range("H37").Value = 42089.176
$ws.Range("J37").Value = 42089.176
$ws.Range("L37").Value = 126267.528
$ws.Range("N37").Value = -126491.528
$ws.Range("H64").Value = 17104
$ws.Range("I64").Value = 25156
$ws.Range("J64").Value = 1000
$ws.Range("K64").Value = 75468
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -75198
$ws.Range("N64").Value = -3540
$ws.Range("H67").Value = 17104
$ws.Range("I67").Value = 25156
$ws.Range("J67").Value = 1000
$ws.Range("K67").Value = 75468
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -74532
$ws.Range("N67").Value = -4872

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10016
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 14666.667
$ws.Range("H80").Value = 20166
$ws.Range("I80").Value = 13999
$ws.Range("K80").Value = 13999
$ws.Range("M80").Value = -13001
$ws.Range("H83").Value = 20166
$ws.Range("I83").Value = 13999
$ws.Range("K83").Value = 69995
$ws.Range("M83").Value = -65003
$ws.Range("H132").Value = 4120.25
$ws.Range("I132").Value = 3994.5715
$ws.Range("K132").Value = 11983.7145
$ws.Range("M132").Value = -9453.7145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2779.8096
$ws.Range("I7").Value = 3055.1428
$ws.Range("J7").Value = 2229.1428
$ws.Range("K7").Value = 3055.1428
$ws.Range("L7").Value = 2229.1428
$ws.Range("M7").Value = -2943.1428
$ws.Range("N7").Value = -2453.1428
$ws.Range("H46").Value = 3378.611
$ws.Range("J46").Value = 3638.7856
$ws.Range("L46").Value = 3638.7856
$ws.Range("N46").Value = -4014.7856
$ws.Range("H61").Value = 965.2857
$ws.Range("I61").Value = 654.9231
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 654.9231
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -452.9231
$ws.Range("N61").Value = -5404
$ws.Range("H113").Value = 965.2857
$ws.Range("I113").Value = 654.9231
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 654.9231
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 1515.0769
$ws.Range("N113").Value = -9340
$ws.Range("H126").Value = 2779.8096
$ws.Range("I126").Value = 3055.1428
$ws.Range("J126").Value = 2229.1428
$ws.Range("K126").Value = 9165.428400000001
$ws.Range("L126").Value = 6687.428400000001
$ws.Range("M126").Value = -6695.428400000001
$ws.Range("N126").Value = -11627.4284

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 27225
$ws.Range("I52").Value = 4450
$ws.Range("K52").Value = 4450
$ws.Range("M52").Value = -4224
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 895.1892
$ws.Range("I132").Value = 746.34283
$ws.Range("K132").Value = 2239.02849
$ws.Range("M132").Value = 290.9715099999999
$ws.Range("H136").Value = 37300.8
$ws.Range("I136").Value = 53502.4
$ws.Range("K136").Value = 160507.2
$ws.Range("M136").Value = -157957.2
